# Update the "保險" (insurance) sheet (sheet5): fill out the full record
# (company, name, owner, property_category, category, date,
# legislator_name, legislator_id, source_file, index) to match the
# pattern used on the sibling "存款" (deposit) sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("保險")
$wsDeposit = $wb.Worksheets.Item("存款")

# --- Extend the header row / data rows with the same cell style as the
#     existing neighbouring column (D) so no new cellXfs entries appear.
$ws.Range("D1").Copy()
$ws.Range("E1:K1").PasteSpecial(-4122)

$ws.Range("D2").Copy()
$ws.Range("E2:K2").PasteSpecial(-4122)

$ws.Range("D3").Copy()
$ws.Range("E3:K3").PasteSpecial(-4122)

# --- Header row (B1:K1): field names ---
$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "property_category"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "legislator_name"
$ws.Range("I1").Value = "legislator_id"
$ws.Range("J1").Value = "source_file"
$ws.Range("K1").Value = "index"

# --- Row 2 (record 92) ---
$ws.Range("B2").Value = "中華郵政"
$ws.Range("C2").Value = "六六金順"
$ws.Range("D2").Value = "楊瓊瓔"
$ws.Range("E2").Value = "insurance"
$ws.Range("F2").Value = "normal"
# Paste the literal date text (avoids Excel auto-converting the string
# "2012-04-26" into a date serial number): copy it as-is from the
# equivalent column on the "存款" sheet, which already holds it as text.
$wsDeposit.Range("I2").Copy()
$ws.Range("G2").PasteSpecial(-4163)
$ws.Range("H2").Value = "楊瓊瓔"
$ws.Range("I2").Value = 854
$ws.Range("J2").Value = "tmp8a701"
$ws.Range("K2").Value = 92

# --- Row 3 (record 93) ---
$ws.Range("B3").Value = "中華郵政"
$ws.Range("C3").Value = "吉利保險"
$ws.Range("D3").Value = "楊瓊瓔"
$ws.Range("E3").Value = "insurance"
$ws.Range("F3").Value = "normal"
$wsDeposit.Range("I2").Copy()
$ws.Range("G3").PasteSpecial(-4163)
$ws.Range("H3").Value = "楊瓊瓔"
$ws.Range("I3").Value = 854
$ws.Range("J3").Value = "tmp8a701"
$ws.Range("K3").Value = 93
